$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Carlos / 004221472 / 1841.84" row (currently row 4).
$ws.Range("A4").EntireRow.Delete() | Out-Null

# Insert two blank rows before the current row 2 so the single
# "008420189 / Styphany / 25000" row can be replaced by three new rows.
$ws.Range("A2:A3").EntireRow.Insert() | Out-Null

# Row 2: new first record (replaces the old Styphany row)
$ws.Range("A2").Value = "'004268684"
$ws.Range("B2").Value = "Patricia"
$ws.Range("C2").Value = 9927.98

# Row 3: new second record
$ws.Range("A3").Value = "'004276856"
$ws.Range("B3").Value = "Daura"
$ws.Range("C3").Value = 9913.98

# Row 4: new third record
$ws.Range("A4").Value = "'004211922"
$ws.Range("B4").Value = "Carlos"
$ws.Range("C4").Value = 3967.63
